$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A86").Value = 85
$ws.Range("B86").Value = "Oxfam Novip (Netherlands)"
$ws.Range("C86").Value = "NGO"

$ws.Range("A87").Value = 86
$ws.Range("B87").Value = "Oxfam Uganda"
$ws.Range("C87").Value = "NGO"
